$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 9) with the same values/shape as row 8,
# matching the new record appended in the diff.
$ws.Range("A9").Value = "ISSFA - 0045"
$ws.Range("B9").Value2 = $ws.Range("B8").Value2
$ws.Range("B9").NumberFormat = $ws.Range("B8").NumberFormat
$ws.Range("C9").Value = "Iñaquito"
$ws.Range("D9").Value = "Quito"
$ws.Range("E9").Value = "Quito"
$ws.Range("F9").Value = "Quito"
$ws.Range("G9").Value = "Pichincha"
$ws.Range("H9").Value = "Casa"
$ws.Range("I9").Value = "Horizontal"
$ws.Range("J9").Value = 523
$ws.Range("K9").Value = 834
$ws.Range("L9").Value = 750.6
$ws.Range("M9").Value = 834
